$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 188, shifting existing rows 188-198 down to 189-199
$ws.Rows.Item(188).Insert()

# Populate the newly inserted row 188 with the new data point
$ws.Cells.Item(188, 1).Value = 7
$ws.Cells.Item(188, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(188, 3).Value = "Ñuble"
$ws.Cells.Item(188, 4).Value = 44615
$ws.Cells.Item(188, 5).Value = 16
$ws.Cells.Item(188, 6).Value = 100112003
$ws.Cells.Item(188, 7).Value = "Ajo"
$ws.Cells.Item(188, 8).Value = "Chino"
$ws.Cells.Item(188, 9).Value = "Primera"
$ws.Cells.Item(188, 10).Value = 60
$ws.Cells.Item(188, 11).Value = 19000
$ws.Cells.Item(188, 12).Value = 20000
$ws.Cells.Item(188, 13).Value = 19500
$ws.Cells.Item(188, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(188, 15).Value = "China"
$ws.Cells.Item(188, 16).Value = 1950
$ws.Cells.Item(188, 17).Value = 10
$ws.Cells.Item(188, 18).Value = "Hortaliza"
